# (XF-876) AUTO_TC 7.3.3 Update User Data - UpdateUser Method - Scroll Methods -
# getListOfElement Method - Getting data from the excel.xlsx
#
# Changes applied to the "2_UserManagement" parameter sheet:
#   - Rename the two test users so they are clearly test data ("Gerardo" ->
#     "tGerardo", "Julio" -> "tJulio") in the "Add FirstName" column.
#   - Add a new "Job Title Change" column (N) used by the new test case, with
#     a "Test Job Title" value for both data rows.
#   - Leave the scroll/selection on the sheet near the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2_UserManagement")

# Rename the FirstName test values.
$ws.Range("E2").Value = "tGerardo"
$ws.Range("E3").Value = "tJulio"

# New "Job Title Change" column (header + the two data rows).
$ws.Range("N1").Value = "Job Title Change"
$ws.Range("N2").Value = "Test Job Title"
$ws.Range("N3").Value = "Test Job Title"

# Update the view: make "2_UserManagement" the active sheet and move the
# selection over to the newly added column.
$ws.Activate()
$ws.Range("M6").Select()
